$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InventoryList")

# ---------------------------------------------------------------------------
# 1. Update the sample inventory rows (4-14) with new Id/Item/Supplier/Group
#    test data (mirrors a manual re-entry of the sample rows in the sheet).
# ---------------------------------------------------------------------------

# Row 4
$ws.Range("C4").Value = "q3"
$ws.Range("D4").Value = "q3"
$ws.Range("E4").Value = "Supik2"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 35
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 10
$ws.Range("M4").Value = "Gr11!"

# Row 5
$ws.Range("C5").Value = "q11"
$ws.Range("D5").Value = "q11"
$ws.Range("E5").Value = "Sup3"
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = 0
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").ClearContents()
$ws.Range("M5").Value = "Gr11!"

# Row 6
$ws.Range("C6").Value = "Id0"
$ws.Range("D6").Value = "item0"
$ws.Range("E6").Value = "Sup1"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").ClearContents()
$ws.Range("M6").Value = "Gr1"

# Row 7
$ws.Range("C7").Value = "id1"
$ws.Range("D7").Value = "item1"
$ws.Range("E7").Value = "Sup2"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").ClearContents()
$ws.Range("M7").Value = "Gr1"

# Row 8
$ws.Range("C8").Value = "id2"
$ws.Range("D8").Value = "item2"
$ws.Range("E8").Value = "Sup3"
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1
$ws.Range("K8").ClearContents()
$ws.Range("M8").Value = "Gr1"

# Row 9
$ws.Range("C9").Value = "id3"
$ws.Range("D9").Value = "item3"
$ws.Range("E9").Value = "Sup1"
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("K9").ClearContents()
$ws.Range("M9").Value = "Gr2"

# Row 10
$ws.Range("C10").Value = "id4"
$ws.Range("D10").Value = "item4"
$ws.Range("E10").Value = "Sup2"
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("K10").ClearContents()
$ws.Range("M10").Value = "Gr2"

# Row 11
$ws.Range("C11").Value = "id5"
$ws.Range("D11").Value = "item5"
$ws.Range("E11").Value = "Sup3"
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 1
$ws.Range("K11").ClearContents()
$ws.Range("M11").Value = "Gr2"

# Row 12
$ws.Range("C12").Value = "id6"
$ws.Range("D12").Value = "item6"
$ws.Range("E12").Value = "Sup1"
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 1
$ws.Range("K12").ClearContents()
$ws.Range("M12").Value = "Gr3"

# Row 13
$ws.Range("C13").Value = "id7"
$ws.Range("D13").Value = "item7"
$ws.Range("E13").Value = "Sup2"
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 1
$ws.Range("K13").ClearContents()
$ws.Range("M13").Value = "Gr3"

# Row 14
$ws.Range("C14").Value = "id8"
$ws.Range("D14").Value = "item8"
$ws.Range("E14").Value = "Sup3"
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 1
$ws.Range("K14").ClearContents()
$ws.Range("M14").Value = "Gr3"

# ---------------------------------------------------------------------------
# 2. Remove 4 blank trailing rows (86:89) shifting the totals row up so the
#    sheet's used range becomes B1:M88.
# ---------------------------------------------------------------------------
$ws.Rows("86:89").Delete()

# Restore the _xlfn.SINGLE(...) formulas on the two rows that shifted up into
# the 86/87 positions (the row-shift rewrite drops the _xlfn. prefix).
$ws.Range("B86").Formula = '=_xlfn.SINGLE(IFERROR((_xlfn.SINGLE(InventoryList!G86:G86)<=_xlfn.SINGLE(InventoryList!I86:I86))*(_xlfn.SINGLE(InventoryList!L86:L86)="")*_xlfn.SINGLE(valHighlight),0))'
$ws.Range("B87").Formula = '=_xlfn.SINGLE(IFERROR((_xlfn.SINGLE(InventoryList!G87:G87)<=_xlfn.SINGLE(InventoryList!I87:I87))*(_xlfn.SINGLE(InventoryList!L87:L87)="")*_xlfn.SINGLE(valHighlight),0))'

# ---------------------------------------------------------------------------
# 3. Selection moves to E4.
# ---------------------------------------------------------------------------
$ws.Range("E4").Select()

$wb.Application.Calculate()
